$d = $word.ActiveDocument

# 1) Rename the first heading: "Objet resourceDetails" -> "Objet geolocalisation"
$d.Content.Find.Execute("Objet resourceDetails", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Objet geolocalisation", 2) | Out-Null

# 2) Insert a new empty Heading1-styled paragraph right after that heading paragraph;
#    it will act as the anchor we replace with the new table + "Type resource" heading.
$p1 = $d.Paragraphs(1)
$p1.Range.InsertParagraphAfter()

# 3) Replace that empty anchor paragraph's XML with the exact target WordprocessingML:
#    the new "resourceDetails" summary table followed by the new "Type resource" heading.
$p2 = $d.Paragraphs(2)
$fragment = @'
<w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:tblPr><w:tblStyle w:val="MediumShading1-Accent1"/><w:tblW w:type="auto" w:w="0"/><w:tblLayout w:type="fixed"/><w:tblLook w:firstColumn="1" w:firstRow="1" w:lastColumn="0" w:lastRow="0" w:noHBand="0" w:noVBand="1" w:val="04A0"/></w:tblPr><w:tblGrid><w:gridCol w:w="2040"/><w:gridCol w:w="2040"/><w:gridCol w:w="2040"/><w:gridCol w:w="2040"/><w:gridCol w:w="2040"/><w:gridCol w:w="2040"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="1701"/></w:tcPr><w:p><w:r><w:t>Nom de balise</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="1984"/></w:tcPr><w:p><w:r><w:t>Champ correspondant</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="1134"/></w:tcPr><w:p><w:r><w:t>Format</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="1417"/></w:tcPr><w:p><w:r><w:t>Cardinalité</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4535"/></w:tcPr><w:p><w:r><w:t>Description</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="1701"/></w:tcPr><w:p><w:r><w:t>Exemple</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="1701"/></w:tcPr><w:p><w:r><w:t>resourceDetails</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="1984"/></w:tcPr><w:p><w:r><w:t>Ressource</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="1134"/></w:tcPr><w:p><w:r><w:t>cf. type resource</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="1417"/></w:tcPr><w:p><w:r><w:t>0..n</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4535"/></w:tcPr><w:p><w:r><w:t>Une liste d’objets Resource détaillant les ressources demandées ainsi que celles notifiées non encore décrites au demandeur</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="1701"/></w:tcPr><w:p><w:r/></w:p></w:tc></w:tr></w:tbl><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Type resource</w:t></w:r></w:p>
'@
$p2.Range.InsertXML($fragment) | Out-Null
